$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new second row of collection data ---
$ws.Range("A2").Value = "MCH243"
$ws.Range("C2").Value = "SCRAP BOOK OF PHOTOGRAPHS & PRESS CUTTINGSS"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"
$ws.Range("H2").Value = ""

# --- Apply the row's font formatting (10pt Calibri, automatic/theme text colour) ---
# Build the font once on a scratch cell, then stamp the resulting format onto each
# cell in the row (keeps the style table compact instead of re-deriving it 7 times).
$scratch = $ws.Range("Z99")
$scratch.Font.Name = "Calibri"

foreach ($addr in @("A2","C2","D2","E2","F2","G2","H2")) {
    $cell = $ws.Range($addr)
    $scratch.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Font.ThemeColor = 1
}
$scratch.Clear()
$excel.CutCopyMode = $false

# --- Re-establish the frozen header pane and select the new row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:H2").Select()
